# Apply the "continue moving flags to most logical position and refactoring code" edit.
#
# Summary of the semantic change (from the OOXML diff):
#   - Row 4 (ht 30 -> 45):
#       * New cell B4 = "mSPCL_EFCT_CONTINUOUS|mSPCL_HANDLER | mSPCL_HANDLER_SHOOT"
#         (styled like B16/B21: vertically centered + wrap text)
#       * D4/E4 changed from "mSPCL_EFCT_CONTINUOUS|mEFCT_UNIQ_WAITING" to "mEFCT_UNIQ_WAITING"
#   - Row 5 (ht 30 -> default/auto):
#       * B5 (was blank, style s="9" vertical-center/no-wrap) = "mSPCL_EFCT_CONTINUOUS"
#         (style is unchanged, only the value is added)
#       * C5 "mROW_MENU" -> blank
#       * D5/E5 changed from "mSPCL_EFCT_CONTINUOUS|mEFCT_UNIQ_WAITING" to "mEFCT_UNIQ_WAITING"
#   - Row 6 (ht 30 -> default/auto):
#       * B6 (was blank, style s="9" vertical-center/no-wrap) = "mSPCL_EFCT_CONTINUOUS"
#         (style is unchanged, only the value is added)
#       * C6 "mROW_MENU" -> blank
#       * D6/E6 changed from "mSPCL_EFCT_CONTINUOUS|mEFCT_UNIQ_WAITING" to "mEFCT_UNIQ_WAITING"
#   - Selection cursor moved from D21:E21 to D16 (cosmetic, applied for completeness)
#
# All other cell-value index changes visible in the raw XML diff are purely due to the
# shared-string table growing by three new entries (shifting later indices) - the actual
# displayed text for those cells is unchanged, so nothing further needs to be done for them;
# Excel/COM manages the shared string table automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 ---
$ws.Range("B4").Value = "mSPCL_EFCT_CONTINUOUS|mSPCL_HANDLER | mSPCL_HANDLER_SHOOT"
$ws.Range("B4").WrapText = $true
$ws.Range("B4").VerticalAlignment = -4108   # xlCenter

$ws.Range("D4").Value = "mEFCT_UNIQ_WAITING"
$ws.Range("E4").Value = "mEFCT_UNIQ_WAITING"

$ws.Rows.Item(4).RowHeight = 45

# --- Row 5 ---
# B5 keeps its existing style (vertical-center, no wrap) - only the value changes.
$ws.Range("B5").Value = "mSPCL_EFCT_CONTINUOUS"

$ws.Range("C5").Value = ""

$ws.Range("D5").Value = "mEFCT_UNIQ_WAITING"
$ws.Range("E5").Value = "mEFCT_UNIQ_WAITING"

$ws.Rows.Item(5).AutoFit()

# --- Row 6 ---
# B6 keeps its existing style (vertical-center, no wrap) - only the value changes.
$ws.Range("B6").Value = "mSPCL_EFCT_CONTINUOUS"

$ws.Range("C6").Value = ""

$ws.Range("D6").Value = "mEFCT_UNIQ_WAITING"
$ws.Range("E6").Value = "mEFCT_UNIQ_WAITING"

$ws.Rows.Item(6).AutoFit()

# --- Selection cursor (cosmetic) ---
$ws.Range("D16").Select()
